# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    2  = 24
    3  = 28
    4  = 6167
    5  = 169
    6  = 21
    8  = 1855
    9  = 1388
    10 = 291
    11 = 947
    12 = 199
    13 = 5560
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
